# Insert a new data row right before the current row 30 ("A1:R86" -> "A1:R87"),
# shifting all existing rows 30..86 down to 31..87, and populate the new
# row 30 with a fresh "Arveja Verde" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 30..86 down by one row.
$ws.Rows.Item(30).Insert()

# Most columns are constant across every data row on this sheet - copy them
# straight from the (still untouched) row immediately above the insertion.
$ws.Range("A30").Value = $ws.Range("A29").Value()
$ws.Range("B30").Value = $ws.Range("B29").Value()
$ws.Range("C30").Value = $ws.Range("C29").Value()
$ws.Range("E30").Value = $ws.Range("E29").Value()
$ws.Range("F30").Value = $ws.Range("F29").Value()
$ws.Range("G30").Value = $ws.Range("G29").Value()
$ws.Range("I30").Value = $ws.Range("I29").Value()
$ws.Range("Q30").Value = $ws.Range("Q29").Value()
$ws.Range("R30").Value = $ws.Range("R29").Value()

# New row-specific values.
$ws.Range("D30").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("D30").Value = 44533
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("J30").Value = 61
$ws.Range("K30").Value = 16000
$ws.Range("L30").Value = 17000
$ws.Range("M30").Value = 16508
$ws.Range("N30").Value = "$/saco 25 kilos"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 660
